# Bug fix in Eduati data files:
#  - Sheet1 had stray leftover rows (45:87) that only carried an
#    incrementing index in column A with no real data - remove them so
#    the sheet matches the real data range (A1:N44), same as Sheet2/Sheet3.
#  - Re-point the "current" sheet/selection at Sheet1 (where the fix was
#    made) instead of Sheet3.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")

# Drop the bogus trailing rows (45-87) left over in column A.
$ws1.Rows("45:87").Delete() | Out-Null

# Make Sheet1 the active sheet/tab, scrolled/selected the way the fixed
# workbook was left (selection on F61, view scrolled so row 24 is at the
# top).
$ws1.Activate() | Out-Null
$ws1.Range("F61").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 24
$excel.ActiveWindow.ScrollColumn = 1
